# presence2.xlsx — "worker and s. manager can get presence report"
#
# The underlying data fix: rows 7 & 8 had been mistakenly entered with new
# employees ("jeck"/"ka", "joni"/"j") instead of reusing the existing
# worker "yoni machluf" for those extra presence entries; row 6's arrival
# date/day was off by one day (16th -> 17th); and the "total seconds"
# column (I) was left blank for the first two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix "total seconds" (column I) left empty on the first two rows ---
$ws.Range("I2").Value = 373
$ws.Range("I3").Value = 112

# --- row 6: arrival was actually the next day (Mon 17th, not Sun 16th) ---
$ws.Range("D6").Value = "Sun, 17 Dec 2018 21:42:40"
$ws.Range("F6").Value = 17

# --- rows 7 & 8: these are further presence entries for the existing
#     worker/s.manager "yoni machluf", not new people ---
$ws.Range("B7").Value = "yoni"
$ws.Range("C7").Value = "machluf"
$ws.Range("B8").Value = "yoni"
$ws.Range("C8").Value = "machluf"

# --- default body font: Calibri -> Arial ---
$wb.Styles.Item(1).Font.Name = "Arial"

# --- widen the date/time columns (D, H) so the full timestamps are
#     visible without truncation ---
$ws.Columns.Item(4).ColumnWidth = 22.75
$ws.Columns.Item(8).ColumnWidth = 22.75

# --- leave the cursor where the author left it ---
$ws.Range("H16").Select() | Out-Null
